$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: Invalid (G) and Absent (H) set to 1
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: Absent (H) set to 1
$ws.Range("H7").Value = 1

# Row 8: Absent (H) set to 1
$ws.Range("H8").Value = 1

# Row 9: Absent (H) set to 1
$ws.Range("H9").Value = 1

# Row 10: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: Absent (H) set to 1
$ws.Range("H11").Value = 1

# Row 12: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: Total Attendance Count (D) and Real (E) set to 1
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 1

# Row 15: Absent (H) set to 1
$ws.Range("H15").Value = 1

# Row 16: Absent (H) set to 1
$ws.Range("H16").Value = 1

# Row 17: Absent (H) set to 1
$ws.Range("H17").Value = 1

# Row 18: Absent (H) set to 1
$ws.Range("H18").Value = 1
